# Began adding support for multitasks
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of the first task's columns (D:H) onto the new
# second-task columns (I:M) so the new cells pick up the same number
# formats / fills as their D:H counterparts, then fill in the values.
$ws.Range("D2:H3").Copy() | Out-Null
$ws.Range("I2:M3").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

# Row 2 - second task values
$ws.Range("I2").Value = 120
$ws.Range("J2").Value = 6.8
$ws.Range("K2").Value = 8
$ws.Range("L2").Value = 300
$ws.Range("M2").Value = 2

# Row 3 - second task values
$ws.Range("I3").Value = 120
$ws.Range("J3").Value = 6.8
$ws.Range("K3").Value = 8
$ws.Range("L3").Value = 10
$ws.Range("M3").Value = 2

# Move the active selection to K4 (matches the recorded cursor position
# after the edit).
$ws.Range("K4").Select() | Out-Null
